$wb = $excel.ActiveWorkbook
$old = $wb.Worksheets.Item(1)
$new = $wb.Worksheets.Add($null, $old)
$new.Range("A1").Value = "Ar da sala"
$new.Range("B1").Value = $old.Range("B1").Value2
$new.Range("C1").Value = $old.Range("C1").Value2
$new.Range("D1").Value = $old.Range("D1").Value2
$new.Name = "Sheet1New"
$new.Outline.SummaryBelow = $false
$new.Activate()
$excel.DisplayAlerts = $false
$old.Delete()
$wb.Worksheets.Item(1).Name = "Sheet1"
$wb.Worksheets.Item(1).Activate()
